$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 12173.1015719068
$ws.Range("D2").Value = 21315.86
$ws.Range("F2").Value = 14.6514194239962

$ws.Range("C3").Value = 8594.84912745595
$ws.Range("F3").Value = 221.293893933667

$ws.Range("C4").Value = 8881.24748000531
$ws.Range("F4").Value = 232.460454083816

$ws.Range("C5").Value = 13354.701608891
$ws.Range("F5").Value = 435.012605476856

$ws.Range("C6").Value = 12770.4982658303
$ws.Range("F6").Value = 398.52377280067

$ws.Range("C7").Value = 12257.3262172698
$ws.Range("F7").Value = 377.087676503394
